$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TABLE_1")
$ws.Range("BH6").Copy()
$ws.Range("BI5:BI56").PasteSpecial(-4122)
$ws.Range("BI5").Value = 21615055.070713539
$ws.Range("BI6").Value = 232687.33251337
$ws.Cells.Columns.Item(61).AutoFit()
